$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.547787547111511
$ws.Range("B1").Value = 3.671063899993896
$ws.Range("C1").Value = 5.683390140533447
$ws.Range("D1").Value = 1.399848341941833
$ws.Range("E1").Value = 0.8167507648468018
